$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SubProcessID"
$ws.Range("B2").Value = "str"

$ws.Range("A3").Value = "concept:name"
$ws.Range("B3").Value = "str"

$ws.Range("A4").Value = "time:timestamp"
$ws.Range("B4").Value = "datetime"

$ws.Range("A5").Value = "stream:datastream"
$ws.Range("B5").Value = "dict"

$ws.Range("A7").Value = "operation_end_time"
$ws.Range("B7").Value = "datetime"
